$wb = $excel.ActiveWorkbook

# The workbook currently has a single sheet "ODI Batting". We need to:
#  1. Update a couple of its values (MATCH_CARD_LINK -> MATCH_CODE header/value)
#  2. Insert a brand new "Player Info" sheet in front of it containing player bio data

$odi = $wb.Worksheets.Item("ODI Batting")
$odi.Range("D1").Value = "MATCH_CODE"
$odi.Range("D2").NumberFormat = "@"
$odi.Range("D2").Value = "4656"

# Create the new sheet and move it so it becomes the first tab in the workbook
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"
$playerInfo.Move($wb.Worksheets.Item(1))

# Header row (bold, centered, bordered - matching the style used on the other sheet)
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$header = $playerInfo.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Data row
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4931"
$playerInfo.Range("B2").Value = "Ruturaj Dashrat Gaikwad"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"
